$d = $word.ActiveDocument
$bm = $d.Bookmarks.Item("_GoBack")
$r = $bm.Range
$delta = 1749 - 427
$r.Move(1, -$delta) | Out-Null
Write-Output ("r Start=" + $r.Start + " End=" + $r.End)
Write-Output ("bm Start=" + $bm.Start + " End=" + $bm.End)
